$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.470.89"
$ws.Range("D3").Value = "2.275.68"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'251.30"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "'0.623"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("D7").Value = "'71.66"
$ws.Range("E7").Value = "  +6.03%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.645"
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("D10").Value = "'38.61"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'59.14"
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "'0.0956"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").Value = "'7.27"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "'0.106"
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").Value = "2.616.54"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "'14.88"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "'0.871"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "2.283.45"
$ws.Range("E18").Value = "  +3.08%  "
$ws.Range("D19").Value = "42.410.84"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").Value = "0.0₃0991"
$ws.Range("E20").Value = "  +2.85%  "
$ws.Range("D21").Value = "'6.29"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").Value = "'71.90"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "'232.39"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "'2.20"
$ws.Range("E24").Value = "  +7.40%  "
$ws.Range("D25").Value = "'3.91"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'11.42"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("D30").Value = "'2.15"
$ws.Range("E30").Value = "  -4.12%  "
$ws.Range("D31").Value = "'167.25"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("D33").Value = "'6.36"
$ws.Range("E33").Value = "  +6.44%  "
$ws.Range("D34").Value = "'0.125"
$ws.Range("E34").Value = "  +3.79%  "
$ws.Range("D35").Value = "'0.0810"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").Value = "'31.13"
$ws.Range("E36").Value = "  +21.20%  "
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("D38").Value = "'4.75"
$ws.Range("E38").Value = "  +15.12%  "
$ws.Range("D39").Value = "'4.74"
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("D40").Value = "'0.0306"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "'13.81"
$ws.Range("E41").Value = "  +13.03%  "
$ws.Range("E42").Value = "  +4.00%  "
$ws.Range("E43").Value = "  +3.90%  "
$ws.Range("D44").Value = "'0.212"
$ws.Range("E44").Value = "  +6.55%  "
$ws.Range("E45").Value = "  +6.54%  "
$ws.Range("D46").Value = "'61.43"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "'4.87"
$ws.Range("E47").Value = "  -4.66%  "
$ws.Range("D48").Value = "'0.104"
$ws.Range("E48").Value = "  +3.33%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "'1.18"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("D51").Value = "'97.39"
$ws.Range("E51").Value = "  +4.34%  "
